$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c1 = $ws.Cells.Item(23,1)
$c1.Value = 45587.669074074074
$c1.NumberFormat = "m/d/yy h:mm"
$c1.Font.Name = "Roboto"
$c1.Font.Color = 4408131

$c2 = $ws.Cells.Item(23,2)
$c2.Value = "The S+ Platform"
$c2.Font.Name = "Roboto"
$c2.Font.Color = 4408131

$blank1 = $ws.Cells.Item(23,28)
$blank1.Font.Name = "Arial"
$blank2 = $ws.Cells.Item(23,29)
$blank2.Font.Name = "Arial"
